$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
